$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($r, $name, $b, $c, $d, $e, $f, $g, $h)
    $ws.Range("A$r").Value = $name
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

# -----------------------------------------------------------------
# Re-order a handful of countries (this also shifts the numeric data
# for the rows that sit between the old and new position by one slot)
# and refresh the COVID figures for a handful of countries.
# -----------------------------------------------------------------

# --- Armenia block (rows 59-61): Armenia moves up, Moldavia/Austria shift down ---
Set-Row 59 "Armenia"  68530 1836 49219 18190 0 20  1121
Set-Row 60 "Moldavia" 67958 0    49083 17258 0 0   1617
Set-Row 61 "Austria"  67451 0    51407 15130 0 0   914

# --- Georgia block (rows 93-94): Georgia moves up, Costa de Marfil shifts down ---
Set-Row 93 "Georgia"          21208 1351 9003  12033 0 14 172
Set-Row 94 "Costa de Marfil"  20342 0    20044 177   0 0  121

# --- Lituania block (rows 116-118): Lituania moves up, Zimbabue/Angola shift down ---
Set-Row 116 "Lituania" 8239 311 3599 4520 0 2 120
Set-Row 117 "Zimbabue" 8187 0   7692 262  0 0 233
Set-Row 118 "Angola"   8049 0   3037 4761 0 0 251

# --- Estonia block (rows 143-144): Estonia moves up, Mayotte shifts down ---
Set-Row 143 "Estonia" 4171 44 3334 766  0 0 71
Set-Row 144 "Mayotte" 4159 0  2964 1152 0 0 43

# --- Fiyi block (rows 205-206): Fiyi moves up, Guam shifts down ---
Set-Row 205 "Fiyi" 33 1 30 1  0 0 2
Set-Row 206 "Guam"  32 0 0  31 0 0 1

# --- Daily figure refresh (no re-ordering involved) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 8520822
$ws.Range("C4").Value = 515
$ws.Range("D4").Value = 5546675
$ws.Range("E4").Value = 2747978
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 226169

# Rusia (row 7)
$ws.Range("B7").Value = 1447335
$ws.Range("C7").Value = 15700
$ws.Range("D7").Value = 1096560
$ws.Range("E7").Value = 325823
$ws.Range("G7").Value = 317
$ws.Range("H7").Value = 24952

# Hungria (row 71)
$ws.Range("B71").Value = 50180
$ws.Range("C71").Value = 1423
$ws.Range("D71").Value = 14905
$ws.Range("E71").Value = 34016
$ws.Range("G71").Value = 48
$ws.Range("H71").Value = 1259
